# Revert "fifth para added":
# Remove the sentence "Preparing food with heat or fire is an activity
# unique to humans." (its own paragraph) and rejoin the remaining text so
# the paragraph that used to read "Cooking can also occur through chemical
# reactions without the presence of heat," flows straight into the
# trailing ". " run, just like before that paragraph was added.

$d = $word.ActiveDocument

# 1. Delete the sentence text itself first, while it still lives in its own
#    paragraph. Doing this before merging the paragraph marks keeps it from
#    being absorbed into neighboring runs.
$find = $d.Content.Find
$find.Text = "Preparing food with heat or fire is an activity unique to humans."
$find.Forward = $true
$find.Wrap = 1
$found = $find.Execute()
if ($found) {
    $find.Parent.Delete()
}

# 2. Delete the paragraph mark that ends the "...without the presence of
#    heat," paragraph, merging it with the (now sentence-less) paragraph
#    that follows. This also drops that paragraph's own pPr, matching the
#    diff.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cooking can also occur through chemical reactions without the presence of heat,*") {
        $target = $p
    }
}
if ($target -ne $null) {
    $r = $target.Range
    $markRange = $d.Range($r.End - 1, $r.End)
    $markRange.Delete()
}
